$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = 8244899
$ws.Range("C32").Value = 3394688
$ws.Range("F32").Value = 3151151
$ws.Range("G32").Value = 28815
$ws.Range("H32").Value = 1269450
$ws.Range("I32").Value = 400795
$ws.Range("J32").Value = 7417789
$ws.Range("K32").Value = 4465835
$ws.Range("L32").Value = 1809324
$ws.Range("M32").Value = 1386
$ws.Range("N32").Value = 1016195
$ws.Range("O32").Value = 68720
$ws.Range("P32").Value = 56328
$ws.Range("Q32").Value = 827111
$ws.Range("R32").Value = 459603
$ws.Range("T32").Value = 454878
$ws.Range("U32").Value = 17897
$ws.Range("W32").Value = 8258071
$ws.Range("X32").Value = 7890564
$ws.Range("Y32").Value = 367507
